$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 7: "co2" node, styled with the grey highlight used for market/balance rows.
$a7 = $ws.Range("A7")
$a7.Interior.Color = 6710886
$a7.Interior.PatternColor = 8421504
$a7.Value = "co2"

# B7 carries the same grey fill but stays empty, like the source row.
$a7.Copy($ws.Range("B7"))
$ws.Range("B7").ClearContents()

# D7 reuses the existing "balance_type_none" label.
$ws.Range("D7").Value = "balance_type_none"

# Move the active selection the way the author's last recorded cursor was.
$null = $ws.Range("C13").Select()
